$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A26").Value = "test"
